$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<ala>"
$ws.Range("C2").Value = 16

# Row 3
$ws.Range("C3").Value = 17

# Row 4
$ws.Range("C4").Value = 14

# Row 5
$ws.Range("C5").Value = 16

# Row 6
$ws.Range("B6").Value = "<see>"

# Row 7
$ws.Range("C7").Value = 12

# Row 8
$ws.Range("C8").Value = 12

# Row 9
$ws.Range("C9").Value = 9

# Row 10
$ws.Range("B10").Value = "<my>"
$ws.Range("C10").Value = 10

# Row 11
$ws.Range("C11").Value = 9

# Row 12
$ws.Range("B12").Value = "<him>"
$ws.Range("C12").Value = 9

# Row 14
$ws.Range("B14").Value = "<all>"
$ws.Range("C14").Value = 9

# Row 16
$ws.Range("C16").Value = 12

# Row 17
$ws.Range("C17").Value = 14

# Row 18
$ws.Range("C18").Value = 10
